# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Reorganiza la tabla de deudores: en vez de agrupar por periodo (con los 5
# trabajadores repetidos por cada periodo), se agrupa por trabajador, listando
# sus 13 periodos de mora en orden descendente (2102 .. 2002). El valor de
# mora es 25749 para el periodo 2102 y 35112 para el resto.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$periods = @("2102","2101","2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002")

$workers = @(
    @("9153983",    "JUAN MANUEL LOMBANA DIAZ"),
    @("1047460015", "LAURA PATRICIA MUÑOZ CARABALLO"),
    @("45555972",   "OLGA PATRICIA POMBO SOTO"),
    @("80874921",   "SAMUEL JOSE POMBO COGOLLO"),
    @("1007130691", "JESUS DAVID PEREZ MARTINEZ")
)

$row = 16
foreach ($worker in $workers) {
    $docNumber = $worker[0]
    $name = $worker[1]
    foreach ($period in $periods) {
        if ($period -eq "2102") {
            $mora = 25749
        } else {
            $mora = 35112
        }
        $ws.Cells.Item($row, 2).Value = "CC"
        $ws.Cells.Item($row, 3).Value = $docNumber
        $ws.Cells.Item($row, 4).Value = $name
        $ws.Cells.Item($row, 5).Value = $period
        $ws.Cells.Item($row, 6).Value = $mora
        $ws.Cells.Item($row, 7).Value = 877803
        $row = $row + 1
    }
}
